$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header tweaks
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 values
$ws.Range("B2").Value = 0.43700851002199509
$ws.Range("C2").Value = 1.9134581246704694
$ws.Range("D2").Value = 1.0116944267165686
$ws.Range("E2").Value = 1.5369554148439599

# Row 3 values
$ws.Range("B3").Value = 1.9833558962570397
$ws.Range("C3").Value = 0.89811331479975587
$ws.Range("D3").Value = 1.3291069666644408
$ws.Range("E3").Value = 1.0594127414566454

# Update selection to match narrowed range used in the diff
$ws.Range("B1:E3").Select()
